# Update Ptdss1-Scarb1 LR-pair sheet with refreshed TPM-derived values
# (ligand/receptor expression + edge weight/specificity columns G,H,I,J,M,N,O,P,Q,R,S,T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.868321
$ws.Range("H2").Value = 11.604963
$ws.Range("I2").Value = 0.1968897496002302
$ws.Range("J2").Value = 0.1968897496002302
$ws.Range("M2").Value = 54.89331066666666
$ws.Range("N2").Value = 164.679932
$ws.Range("O2").Value = 0.8971624892852424
$ws.Range("P2").Value = 0.8971624892852424
$ws.Range("Q2").Value = 212.3449464113906
$ws.Range("R2").Value = 1911.104517702516
$ws.Range("S2").Value = 0.1766420978660906
$ws.Range("T2").Value = 0.1766420978660906
$ws.Range("G3").Value = 3.868321
$ws.Range("H3").Value = 11.604963
$ws.Range("I3").Value = 0.1968897496002302
$ws.Range("J3").Value = 0.1968897496002302
$ws.Range("M3").Value = 2.309992333333333
$ws.Range("N3").Value = 6.929977
$ws.Range("O3").Value = 0.03775393480250816
$ws.Range("P3").Value = 0.03775393480250816
$ws.Range("Q3").Value = 8.935791852872333
$ws.Range("R3").Value = 80.422126675851
$ws.Range("S3").Value = 0.007433362769689249
$ws.Range("T3").Value = 0.007433362769689249
$ws.Range("G4").Value = 3.868321
$ws.Range("H4").Value = 11.604963
$ws.Range("I4").Value = 0.1968897496002302
$ws.Range("J4").Value = 0.1968897496002302
$ws.Range("M4").Value = 3.982169333333333
$ws.Range("N4").Value = 11.946508
$ws.Range("O4").Value = 0.06508357591224938
$ws.Range("P4").Value = 0.06508357591224936
$ws.Range("Q4").Value = 15.40430925768933
$ws.Range("R4").Value = 138.638783319204
$ws.Range("S4").Value = 0.01281428896445036
$ws.Range("T4").Value = 0.01281428896445035
$ws.Range("I5").Value = 0.5409510947037116
$ws.Range("J5").Value = 0.5409510947037116
$ws.Range("M5").Value = 54.89331066666666
$ws.Range("N5").Value = 164.679932
$ws.Range("O5").Value = 0.8971624892852424
$ws.Range("P5").Value = 0.8971624892852424
$ws.Range("Q5").Value = 583.4139738065288
$ws.Range("R5").Value = 5250.725764258759
$ws.Range("S5").Value = 0.4853210307059588
$ws.Range("T5").Value = 0.4853210307059588
$ws.Range("I6").Value = 0.5409510947037116
$ws.Range("J6").Value = 0.5409510947037116
$ws.Range("M6").Value = 2.309992333333333
$ws.Range("N6").Value = 6.929977
$ws.Range("O6").Value = 0.03775393480250816
$ws.Range("P6").Value = 0.03775393480250816
$ws.Range("S6").Value = 0.02042303236078935
$ws.Range("T6").Value = 0.02042303236078935
$ws.Range("I7").Value = 0.5409510947037116
$ws.Range("J7").Value = 0.5409510947037116
$ws.Range("M7").Value = 3.982169333333333
$ws.Range("N7").Value = 11.946508
$ws.Range("O7").Value = 0.06508357591224938
$ws.Range("P7").Value = 0.06508357591224936
$ws.Range("Q7").Value = 42.32306645227111
$ws.Range("R7").Value = 380.90759807044
$ws.Range("S7").Value = 0.03520703163696342
$ws.Range("T7").Value = 0.03520703163696341
$ws.Range("G8").Value = 5.150678333333333
$ws.Range("H8").Value = 15.452035
$ws.Range("I8").Value = 0.2621591556960581
$ws.Range("J8").Value = 0.2621591556960581
$ws.Range("M8").Value = 54.89331066666666
$ws.Range("N8").Value = 164.679932
$ws.Range("O8").Value = 0.8971624892852424
$ws.Range("P8").Value = 0.8971624892852424
$ws.Range("Q8").Value = 282.7377858957354
$ws.Range("R8").Value = 2544.64007306162
$ws.Range("S8").Value = 0.2351993607131929
$ws.Range("T8").Value = 0.2351993607131929
$ws.Range("G9").Value = 5.150678333333333
$ws.Range("H9").Value = 15.452035
$ws.Range("I9").Value = 0.2621591556960581
$ws.Range("J9").Value = 0.2621591556960581
$ws.Range("M9").Value = 2.309992333333333
$ws.Range("N9").Value = 6.929977
$ws.Range("O9").Value = 0.03775393480250816
$ws.Range("P9").Value = 0.03775393480250816
$ws.Range("Q9").Value = 11.89802746146611
$ws.Range("R9").Value = 107.082247153195
$ws.Range("S9").Value = 0.009897539672029563
$ws.Range("T9").Value = 0.009897539672029563
$ws.Range("G10").Value = 5.150678333333333
$ws.Range("H10").Value = 15.452035
$ws.Range("I10").Value = 0.2621591556960581
$ws.Range("J10").Value = 0.2621591556960581
$ws.Range("M10").Value = 3.982169333333333
$ws.Range("N10").Value = 11.946508
$ws.Range("O10").Value = 0.06508357591224938
$ws.Range("P10").Value = 0.06508357591224936
$ws.Range("Q10").Value = 20.51087330486444
$ws.Range("R10").Value = 184.59785974378
$ws.Range("S10").Value = 0.0170622553108356
$ws.Range("T10").Value = 0.0170622553108356
